$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.744.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.884.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.74%  "
$ws.Range("E4").Value = "  +0.66%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("E7").Value = "  -1.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3790"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.55"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07706"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9596"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.887.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.936"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.647"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06999"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.004"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "83.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009471"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.698.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.313"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.112.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.088"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.591"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.64%  "
$ws.Range("E31").Value = "  -5.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09221"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8421"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.056"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.233"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.981"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05651"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.74%  "
$ws.Range("E38").Value = "  -3.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.002"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("E40").Value = "  -4.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.409"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5470"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.80%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1742"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.99%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.000002978"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -30.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.139"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.693"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5141"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.10%  "
$ws.Range("E48").Value = "  -7.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06795"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("E50").Value = "  -5.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.23%  "
